$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.38281187720568
$ws.Range("C2").Value = 6.832164312212846
$ws.Range("D2").Value = 7.896926053294457
$ws.Range("E2").Value = 12.68668211513308
$ws.Range("F2").Value = 37.45241620910736
$ws.Range("I2").Value = 28.21906590238764
$ws.Range("J2").Value = 10.0360696387417
$ws.Range("K2").Value = 10.72411960809379
$ws.Range("L2").Value = 11.0666831581491
$ws.Range("O2").Value = 28.97586259261847
$ws.Range("B3").Value = 13.15373421135708
$ws.Range("C3").Value = 6.775562307236231
$ws.Range("D3").Value = 7.876471785846128
$ws.Range("E3").Value = 12.70393993533293
$ws.Range("F3").Value = 37.55689932773207
$ws.Range("I3").Value = 28.32461286280017
$ws.Range("J3").Value = 10.05583930295011
$ws.Range("K3").Value = 10.56022031431842
$ws.Range("L3").Value = 11.06398152686165
$ws.Range("O3").Value = 29.07451256594319
$ws.Range("B4").Value = 13.01292755549217
$ws.Range("C4").Value = 6.740177643212042
$ws.Range("D4").Value = 7.864896692918507
$ws.Range("E4").Value = 12.71604418821481
$ws.Range("F4").Value = 37.62830855442296
$ws.Range("I4").Value = 28.3939580163021
$ws.Range("J4").Value = 10.06870283130541
$ws.Range("K4").Value = 10.45963443465757
$ws.Range("L4").Value = 11.0636389974305
$ws.Range("O4").Value = 29.14017482780641
$ws.Range("B5").Value = 12.95558083655271
$ws.Range("C5").Value = 6.725605647484109
$ws.Range("D5").Value = 7.860430381992948
$ws.Range("E5").Value = 12.72135647915251
$ws.Range("F5").Value = 37.65923092040789
$ws.Range("I5").Value = 28.42335842002165
$ws.Range("J5").Value = 10.07412755763591
$ws.Range("K5").Value = 10.41870426955324
$ws.Range("L5").Value = 11.06383173903081
$ws.Range("O5").Value = 29.16821213268872
$ws.Range("B6").Value = 12.94606256384324
$ws.Range("C6").Value = 6.723176951551142
$ws.Range("D6").Value = 7.85970398255685
$ws.Range("E6").Value = 12.72226152972425
$ws.Range("F6").Value = 37.66447556738157
$ws.Range("I6").Value = 28.42830930071982
$ws.Range("J6").Value = 10.0750393795069
$ws.Range("K6").Value = 10.41191281736051
$ws.Range("L6").Value = 11.06388385230204
$ws.Range("O6").Value = 29.1729449572438
$ws.Range("B7").Value = 13.01215392839103
$ws.Range("C7").Value = 6.739981729303381
$ws.Range("D7").Value = 7.864835439759271
$ws.Range("E7").Value = 12.71611429350102
$ws.Range("F7").Value = 37.62871820782303
$ws.Range("I7").Value = 28.39434989729133
$ws.Range("J7").Value = 10.06877525056701
$ws.Range("K7").Value = 10.45908213316095
$ws.Range("L7").Value = 11.06364024954631
$ws.Range("O7").Value = 29.14054776942824
$ws.Range("B8").Value = 13.30389992703903
$ws.Range("C8").Value = 6.812782541872305
$ws.Range("D8").Value = 7.88967119748086
$ws.Range("E8").Value = 12.69231998383495
$ws.Range("F8").Value = 37.48693475803132
$ws.Range("I8").Value = 28.25451665616566
$ws.Range("J8").Value = 10.04273604730381
$ws.Range("K8").Value = 10.66762508636218
$ws.Range("L8").Value = 11.06547926912064
$ws.Range("O8").Value = 29.00882008769197
$ws.Range("B9").Value = 13.87173661705169
$ws.Range("C9").Value = 6.950273706438755
$ws.Range("D9").Value = 7.946032235522823
$ws.Range("E9").Value = 12.6575993939333
$ws.Range("F9").Value = 37.26656081265264
$ws.Range("I9").Value = 28.0163022252545
$ws.Range("J9").Value = 9.997404019970977
$ws.Range("K9").Value = 11.07492515923826
$ws.Range("L9").Value = 11.07946813859813
$ws.Range("O9").Value = 28.79092703572797
$ws.Range("B10").Value = 14.28233424825475
$ws.Range("C10").Value = 7.047726702871164
$ws.Range("D10").Value = 7.991899922352068
$ws.Range("E10").Value = 12.63933534555411
$ws.Range("F10").Value = 37.13991799865094
$ws.Range("I10").Value = 27.86320918337104
$ws.Range("J10").Value = 9.967563316600122
$ws.Range("K10").Value = 11.3704901918383
$ws.Range("L10").Value = 11.09598834228605
$ws.Range("O10").Value = 28.65552698979307
$ws.Range("B11").Value = 14.46689614775692
$ws.Range("C11").Value = 7.091214984104675
$ws.Range("D11").Value = 8.013686470472352
$ws.Range("E11").Value = 12.63259190499313
$ws.Range("F11").Value = 37.08998309681505
$ws.Range("I11").Value = 27.7983194297988
$ws.Range("J11").Value = 9.954734229075743
$ws.Range("K11").Value = 11.50360955901335
$ws.Range("L11").Value = 11.10483874285127
$ws.Range("O11").Value = 28.59929916586535
$ws.Range("B12").Value = 14.53640271314378
$ws.Range("C12").Value = 7.107555553865817
$ws.Range("D12").Value = 8.022064575736074
$ws.Range("E12").Value = 12.63026264678011
$ws.Range("F12").Value = 37.07217913818896
$ws.Range("I12").Value = 27.77443069408994
$ws.Range("J12").Value = 9.949982941436881
$ws.Range("K12").Value = 11.55378327987766
$ws.Range("L12").Value = 11.10838022199973
$ws.Range("O12").Value = 28.57877960485483
$ws.Range("B13").Value = 14.52145124534604
$ws.Range("C13").Value = 7.104042087601646
$ws.Range("D13").Value = 8.020254579792905
$ws.Range("E13").Value = 12.63075432845083
$ws.Range("F13").Value = 37.07596435793604
$ws.Range("I13").Value = 27.77954515492718
$ws.Range("J13").Value = 9.95100147229831
$ws.Range("K13").Value = 11.54298863155082
$ws.Range("L13").Value = 11.10760908248927
$ws.Range("O13").Value = 28.5831644775969
$ws.Range("B14").Value = 14.47262252435055
$ws.Range("C14").Value = 7.092561922901425
$ws.Range("D14").Value = 8.014373202648798
$ws.Range("E14").Value = 12.63239578366419
$ws.Range("F14").Value = 37.08849619428572
$ws.Range("I14").Value = 27.79634038784788
$ws.Range("J14").Value = 9.954341199477955
$ws.Range("K14").Value = 11.50774234916415
$ws.Range("L14").Value = 11.10512630585988
$ws.Range("O14").Value = 28.5975955174915
$ws.Range("B15").Value = 14.44266180715483
$ws.Range("C15").Value = 7.085513204794635
$ws.Range("D15").Value = 8.010787225928865
$ws.Range("E15").Value = 12.63343041598844
$ws.Range("F15").Value = 37.09631629193376
$ws.Range("I15").Value = 27.80671698276824
$ws.Range("J15").Value = 9.956400776542681
$ws.Range("K15").Value = 11.48612102088532
$ws.Range("L15").Value = 11.1036302194872
$ws.Range("O15").Value = 28.606535603822
$ws.Range("B16").Value = 14.2702224901372
$ws.Range("C16").Value = 7.044867134276597
$ws.Range("D16").Value = 7.990494278993592
$ws.Range("E16").Value = 12.63980747717008
$ws.Range("F16").Value = 37.14333595302936
$ws.Range("I16").Value = 27.86754553467507
$ws.Range("J16").Value = 9.968416696118155
$ws.Range("K16").Value = 11.36175988882496
$ws.Range("L16").Value = 11.09543664962552
$ws.Range("O16").Value = 28.65930968274081
$ws.Range("B17").Value = 14.16382143514835
$ws.Range("C17").Value = 7.019711581304169
$ws.Range("D17").Value = 7.978278019206965
$ws.Range("E17").Value = 12.64411993560546
$ws.Range("F17").Value = 37.17414807369148
$ws.Range("I17").Value = 27.90607932809586
$ws.Range("J17").Value = 9.97597874842527
$ws.Range("K17").Value = 11.28509481925741
$ws.Range("L17").Value = 11.09075077724923
$ws.Range("O17").Value = 28.69306005390418
$ws.Range("B18").Value = 14.10241692470868
$ws.Range("C18").Value = 7.005163734279654
$ws.Range("D18").Value = 7.971338524419953
$ws.Range("E18").Value = 12.64674769446141
$ws.Range("F18").Value = 37.19259276978358
$ws.Range("I18").Value = 27.92869031617536
$ws.Range("J18").Value = 9.980398444178265
$ws.Range("K18").Value = 11.24087598551054
$ws.Range("L18").Value = 11.08818139574691
$ws.Range("O18").Value = 28.71297741995784
$ws.Range("B19").Value = 14.08159311839189
$ws.Range("C19").Value = 7.00022469905544
$ws.Range("D19").Value = 7.969004001555671
$ws.Range("E19").Value = 12.64766273617582
$ws.Range("F19").Value = 37.19896184474592
$ws.Range("I19").Value = 27.93642284101613
$ws.Range("J19").Value = 9.981906947161553
$ws.Range("K19").Value = 11.22588446629309
$ws.Range("L19").Value = 11.08733311378244
$ws.Range("O19").Value = 28.7198078083261
$ws.Range("B20").Value = 14.17516974204152
$ws.Range("C20").Value = 7.022397658286804
$ws.Range("D20").Value = 7.979569491059833
$ws.Range("E20").Value = 12.64364562217038
$ws.Range("F20").Value = 37.17079329519772
$ws.Range("I20").Value = 27.90193103813716
$ws.Range("J20").Value = 9.975166492231793
$ws.Range("K20").Value = 11.29326901135252
$ws.Range("L20").Value = 11.09123658972174
$ws.Range("O20").Value = 28.68941499148998
$ws.Range("B21").Value = 14.48697558996007
$ws.Range("C21").Value = 7.095937434564789
$ws.Range("D21").Value = 8.016097267604279
$ws.Range("E21").Value = 12.6319075658727
$ws.Range("F21").Value = 37.08478528113565
$ws.Range("I21").Value = 27.79138866682557
$ws.Range("J21").Value = 9.953357345093794
$ws.Range("K21").Value = 11.51810177320666
$ws.Range("L21").Value = 11.10585041669538
$ws.Range("O21").Value = 28.59333579249645
$ws.Range("B22").Value = 14.68849893577606
$ws.Range("C22").Value = 7.143253988765362
$ws.Range("D22").Value = 8.040714271569314
$ws.Range("E22").Value = 12.6255433151684
$ws.Range("F22").Value = 37.03501720018929
$ws.Range("I22").Value = 27.72312749127765
$ws.Range("J22").Value = 9.939726187723027
$ws.Range("K22").Value = 11.66365017771085
$ws.Range("L22").Value = 11.11650800842901
$ws.Range("O22").Value = 28.53504679938766
$ws.Range("B23").Value = 14.58116910890301
$ws.Range("C23").Value = 7.118070519978026
$ws.Range("D23").Value = 8.02750914351757
$ws.Range("E23").Value = 12.62882066575804
$ws.Range("F23").Value = 37.06098935579104
$ws.Range("I23").Value = 27.75919512094097
$ws.Range("J23").Value = 9.946944578342018
$ws.Range("K23").Value = 11.58610968233021
$ws.Range("L23").Value = 11.11071928027726
$ws.Range("O23").Value = 28.56574425533584
$ws.Range("B24").Value = 14.17003989375936
$ws.Range("C24").Value = 7.021183549002555
$ws.Range("D24").Value = 7.978985355297898
$ws.Range("E24").Value = 12.64385959676398
$ws.Range("F24").Value = 37.1723077155341
$ws.Range("I24").Value = 27.90380505526326
$ws.Range("J24").Value = 9.97553348841962
$ws.Range("K24").Value = 11.28957390136887
$ws.Range("L24").Value = 11.09101656567733
$ws.Range("O24").Value = 28.69106132373196
$ws.Range("B25").Value = 13.71899948496748
$ws.Range("C25").Value = 6.913679272198977
$ws.Range("D25").Value = 7.929986127447375
$ws.Range("E25").Value = 12.66571731270729
$ws.Range("F25").Value = 37.31999215324016
$ws.Range("I25").Value = 28.07689411168266
$ws.Range("J25").Value = 10.00905701835516
$ws.Range("K25").Value = 10.96519671820617
$ws.Range("L25").Value = 11.07458079971947
$ws.Range("O25").Value = 28.84554074526551
